# Remove the italicized "Ruth" subtitle paragraph that immediately
# follows the "RUT" Heading2 paragraph at the start of the Ruth intro
# section. (The fuller "Ruth" section later in the document is left
# untouched.)

$d = $word.ActiveDocument

# Locate the run by its unique formatting (italic) rather than a fixed
# paragraph index, since it's the only italic text in the document.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Font.Italic = $true
$found = $rng.Find.Execute("Ruth")

if ($found) {
    $paraIndex = $rng.Paragraphs.Item(1).Index
    $target = $d.Paragraphs.Item($paraIndex)
    # Delete the whole paragraph, including its paragraph mark, so the
    # surrounding paragraphs collapse together exactly as in the diff.
    $target.Range.Delete()
}
